$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction7")

$values = @(1,2,6,7,10,11,12,13,14,15,16,17,19,20,21,22,23,24,25,26,27,29,31,33,34,37,38,39,40,41,42,43,44,45,46,48,49,51,52,53,54,56,59,60,61,63,64,65,66,67,68,69,70,71,72,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
